$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.638.96"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.91%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.393.10"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +0.01%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.01%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.68"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.68%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.89"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  -0.89%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.03%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.972.49"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  +1.66%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.397.22"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  -0.32%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.676.90"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("E18").Value = "  +0.05%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.64"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.95%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.05"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +1.69%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.87"
$ws.Range("D21").Style = $style
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.78"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.27%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -3.36%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.197"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +8.97%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.05%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.63%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +0.33%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("E33").Value = "  -0.53%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.91"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -0.94%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "169.40"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +1.72%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.05"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +1.02%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.425.06"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -0.38%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0766"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -0.50%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.94"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("E41").Value = "  -0.20%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.43"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  +1.35%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.451.85"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.83%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.73"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -1.88%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -0.11%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0263"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("E51").Value = "  -1.00%  "
